$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reactions")
$ws.Columns.Item(6).Insert()
$ws.Cells.Item(1, 6).Value = "Rate units"
$ws.Cells.Item(2, 6).Value = "s^-1"
$ws.AutoFilterMode = $false
$ws.Range("A1:J2").AutoFilter()

$ws2 = $wb.Worksheets.Item("Functions")
$ws2.AutoFilterMode = $false
$ws2.Range("A1:D2").AutoFilter()

$ws3 = $wb.Worksheets.Item("Parameters")
$ws3.AutoFilterMode = $false
$ws3.Range("A1:H7").AutoFilter()
